$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Alkane_Series")
$ws2 = $wb.Worksheets.Item("Compounds_of_interest")

# Rename the "RT_seconds" header to "retention_time" on both sheets
# (they shared the same string value).
$ws1.Range("C1").Value = "retention_time"
$ws2.Range("B1").Value = "retention_time"

# Update the active selection on each sheet.
$ws1.Range("C2").Select()
$ws2.Range("B2").Select()
